$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '67.706.48'
$ws.Range("E2").Value = '  -1.01%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.784.89'
$ws.Range("E3").Value = '  -1.51%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.08%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '596.70'
$ws.Range("E5").Value = '  -0.96%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '169.62'
$ws.Range("E6").Value = '  +0.35%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.782.67'
$ws.Range("E7").Value = '  -1.57%  '

$ws.Range("E8").Value = '  +0.01%  '

$ws.Range("E9").Value = '  -0.77%  '

$ws.Range("E10").Value = '  -1.33%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.48'
$ws.Range("E11").Value = '  -0.13%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.454'
$ws.Range("E12").Value = '  -0.91%  '

$ws.Range("E13").Value = '  +2.09%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '36.62'
$ws.Range("E14").Value = '  -1.32%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.419.85'
$ws.Range("E15").Value = '  -1.47%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.781.89'
$ws.Range("E16").Value = '  -1.28%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '18.69'
$ws.Range("E17").Value = '  +1.12%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '67.726.45'
$ws.Range("E18").Value = '  -1.17%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.20'
$ws.Range("E19").Value = '  -2.76%  '

$ws.Range("E20").Value = '  +0.81%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.60'
$ws.Range("E21").Value = '  -4.39%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '469.03'
$ws.Range("E22").Value = '  -0.34%  '

$ws.Range("E23").Value = '  -2.03%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.0000148'
$ws.Range("E24").Value = '  -8.17%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '83.72'
$ws.Range("E25").Value = '  +0.36%  '

$ws.Range("E26").Value = '  -0.80%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '12.17'
$ws.Range("E27").Value = '  +0.15%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.33'
$ws.Range("E28").Value = '  +2.46%  '

$ws.Range("E29").Value = '  -0.02%  '

$ws.Range("E30").Value = '  -1.87%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.931.21'
$ws.Range("E31").Value = '  -1.52%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.66'
$ws.Range("E32").Value = '  -0.38%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '30.53'
$ws.Range("E33").Value = '  -3.34%  '

$ws.Range("E34").Value = '  -3.51%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '9.13'
$ws.Range("E35").Value = '  -2.69%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.746.93'
$ws.Range("E36").Value = '  -1.57%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.84'
$ws.Range("E37").Value = '  +3.71%  '

$ws.Range("E38").Value = '  -0.99%  '

$ws.Range("B39").Value = 'Kaspa'
$ws.Range("C39").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.138'
$ws.Range("E39").Value = '  -1.18%  '

$ws.Range("B40").Value = 'Mantle'
$ws.Range("C40").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.00'
$ws.Range("E40").Value = '  -2.18%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.80'
$ws.Range("E41").Value = '  -2.39%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.00'
$ws.Range("E42").Value = '  +0.01%  '

$ws.Range("E43").Value = '  -0.88%  '

$ws.Range("E45").Value = '  -0.14%  '

$ws.Range("E46").Value = '  -1.94%  '

$ws.Range("E47").Value = '  -2.77%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '397.42'
$ws.Range("E48").Value = '  -4.38%  '

$ws.Range("E49").Value = '  -6.79%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '139.54'
$ws.Range("E50").Value = '  -1.32%  '

$ws.Range("E51").Value = '  -2.04%  '
